# Weekly time record update for week of 05 September 2016.
# Fills in the actual clock-in/clock-out times for Thursday (row16),
# Friday (row17), Saturday (row18) and Sunday (row19), and clears the
# leftover "1 - " placeholder note that used to sit in column L / row16
# now that the real times have been entered.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Thursday: 1:00 PM - 4:00 PM (3 hrs)
$ws.Range("C16").Value = 0.541666666666667
$ws.Range("D16").Value = 0.666666666666667

# Friday: 10:00 AM - 11:00 AM (1 hr)
$ws.Range("C17").Value = 0.416666666666667
$ws.Range("D17").Value = 0.458333333333333

# Saturday: 6:00 PM - 7:00 PM (1 hr)
$ws.Range("C18").Value = 0.75
$ws.Range("D18").Value = 0.791666666666667

# Sunday: ~9:00 PM - midnight (3 hrs)
$ws.Range("C19").Value = 0.874305555555556
$ws.Range("D19").Value = 0.999305555555556

# The "1 - " placeholder note in L16 is no longer needed now that the
# actual Thursday times are filled in above.
$ws.Range("L16").Value = ""
